$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple F/G quantity+value updates ---
$ws.Range("F23").Value = 11
$ws.Range("G23").Value = 450.89
$ws.Range("F26").Value = 21
$ws.Range("G26").Value = 537.81
$ws.Range("F27").Value = 51
$ws.Range("G27").Value = 1828.86
$ws.Range("F29").Value = 71
$ws.Range("G29").Value = 3637.33
$ws.Range("F114").Value = 44
$ws.Range("G114").Value = 2056.56
$ws.Range("F151").Value = 34
$ws.Range("G151").Value = 4534.24
$ws.Range("F173").Value = 47
$ws.Range("G173").Value = 3694.67
$ws.Range("F178").Value = 104
$ws.Range("G178").Value = 10063.04
$ws.Range("F182").Value = 18
$ws.Range("G182").Value = 1611.72
$ws.Range("F186").Value = 17
$ws.Range("G186").Value = 735.76
$ws.Range("F267").Value = 129
$ws.Range("G267").Value = 5479.92
$ws.Range("F280").Value = 13
$ws.Range("G280").Value = 1262.3
$ws.Range("F370").Value = 219
$ws.Range("G370").Value = 36351.81
$ws.Range("F377").Value = 46
$ws.Range("G377").Value = 44611.26
$ws.Range("F387").Value = 429
$ws.Range("G387").Value = 41441.4
$ws.Range("F399").Value = 267
$ws.Range("G399").Value = 26232.75
$ws.Range("F433").Value = 125
$ws.Range("G433").Value = 1205
$ws.Range("F519").Value = 366
$ws.Range("G519").Value = 20086.08
$ws.Range("F528").Value = 273
$ws.Range("G528").Value = 4329.78
$ws.Range("F558").Value = 198
$ws.Range("G558").Value = 24126.3
$ws.Range("F620").Value = 363
$ws.Range("G620").Value = 28528.17

# --- Subtotal-only updates (no item-level breakdown other than above) ---
$ws.Range("B34").Value = 58828.69
$ws.Range("B123").Value = 73132.67
$ws.Range("B155").Value = 37147.22
$ws.Range("B193").Value = 64670.36
$ws.Range("B295").Value = 117684.17
$ws.Range("B372").Value = 60879.37
$ws.Range("B378").Value = 44611.26
$ws.Range("B389").Value = 58382.86
$ws.Range("B417").Value = 172039.09
$ws.Range("B438").Value = 25025.61
$ws.Range("B525").Value = 120976.67
$ws.Range("B535").Value = 24020.47
$ws.Range("B561").Value = 28616.5
$ws.Range("B628").Value = 210576.57
$ws.Range("B718").Value = 2716703.57
$ws.Range("B719").Value = 2716703.57

# --- Row swaps: 290/291 (code+name+qty+value swap) ---
$ws.Range("B290").Value = 64983
$ws.Range("C290").Value = 'HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S'
$ws.Range("F290").Value = 6
$ws.Range("G290").Value = 514.08
$ws.Range("B291").Value = 66194
$ws.Range("C291").Value = 'HIM-Total Care Baby Pants Diapers-M-9s'
$ws.Range("F291").Value = 23
$ws.Range("G291").Value = 1970.64

# --- Row swaps: 292/293 ---
$ws.Range("B292").Value = 66196
$ws.Range("C292").Value = 'HIM-Total Care Baby Pants Drapers-Xl-9S'
$ws.Range("F292").Value = 6
$ws.Range("G292").Value = 526.2
$ws.Range("B293").Value = 64985
$ws.Range("C293").Value = 'HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S'
$ws.Range("F293").Value = 12
$ws.Range("G293").Value = 1052.4

# --- Row swaps: 297/298 (also D/E) ---
$ws.Range("B297").Value = 61610
$ws.Range("E297").Value = 122.71
$ws.Range("F297").Value = -58
$ws.Range("G297").Value = -5957.18
$ws.Range("B298").Value = 63565
$ws.Range("E298").Value = 109.19
$ws.Range("F298").Value = 60
$ws.Range("G298").Value = 6162.6

# --- Row swaps: 304/305 (also D/E) ---
$ws.Range("B304").Value = 55373
$ws.Range("E304").Value = 163.62
$ws.Range("F304").Value = -94
$ws.Range("G304").Value = -13562.32
$ws.Range("B305").Value = 63520
$ws.Range("E305").Value = 153.4
$ws.Range("F305").Value = 38
$ws.Range("G305").Value = 5482.64

# --- Row swaps: 381/382 (also D/E) ---
$ws.Range("B381").Value = 47097
$ws.Range("D381").Value = 112.28
$ws.Range("E381").Value = 134.16
$ws.Range("F381").Value = 15
$ws.Range("G381").Value = 1684.2
$ws.Range("B382").Value = 58047
$ws.Range("D382").Value = 105.54
$ws.Range("E382").Value = 126.1
$ws.Range("F382").Value = 32
$ws.Range("G382").Value = 3377.28

# --- Row swaps: 479/480 (also E) ---
$ws.Range("B479").Value = 64810
$ws.Range("E479").Value = 291.22
$ws.Range("F479").Value = 0
$ws.Range("G479").Value = 0
$ws.Range("B480").Value = 53319
$ws.Range("E480").Value = 310.64
$ws.Range("F480").Value = -6
$ws.Range("G480").Value = -1643.52
